$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'281.14"
$ws.Range("E2").Value = "'5.87%"
$ws.Range("D3").Value = "'26.83"
$ws.Range("E3").Value = "'0.60%"
$ws.Range("D4").Value = "'4.938"
$ws.Range("E4").Value = "'4.96%"
$ws.Range("E5").Value = "'5.37%"
$ws.Range("D6").Value = "'7.010"
$ws.Range("E6").Value = "'4.06%"
$ws.Range("D7").Value = "'3.354"
$ws.Range("E7").Value = "'5.81%"
$ws.Range("D8").Value = "'0.8855"
$ws.Range("D9").Value = "'1.008"
$ws.Range("E9").Value = "'11.12%"
$ws.Range("D10").Value = "'0.1489"
$ws.Range("E10").Value = "'5.59%"
$ws.Range("D11").Value = "'0.05169"
$ws.Range("E11").Value = "'1.38%"
$ws.Range("D12").Value = "'0.07397"
$ws.Range("E12").Value = "'4.33%"
$ws.Range("D13").Value = "'0.03109"
$ws.Range("E13").Value = "'-1.27%"
$ws.Range("D14").Value = "'0.09056"
$ws.Range("E14").Value = "'0.32%"
$ws.Range("D15").Value = "'0.001560"
$ws.Range("E15").Value = "'1.66%"
$ws.Range("D16").Value = "'0.0006325"
$ws.Range("E16").Value = "'4.44%"
$ws.Range("D17").Value = "'0.006032"
$ws.Range("E17").Value = "'-1.23%"
$ws.Range("D18").Value = "'3.509"
$ws.Range("E18").Value = "'1.52%"
$ws.Range("E19").Value = "'5.66%"
$ws.Range("D20").Value = "'0.3114"
$ws.Range("E20").Value = "'1.40%"
$ws.Range("D21").Value = "'0.1329"
$ws.Range("E21").Value = "'3.75%"
$ws.Range("D22").Value = "'3.942"
$ws.Range("E22").Value = "'-3.47%"
$ws.Range("D23").Value = "'0.04353"
$ws.Range("E23").Value = "'2.44%"
$ws.Range("D24").Value = "'0.001177"
$ws.Range("E24").Value = "'-0.22%"
$ws.Range("D25").Value = "'0.003687"
$ws.Range("E25").Value = "'-9.19%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.15%"
$ws.Range("D27").Value = "'0.0001693"
$ws.Range("E27").Value = "'0.66%"
$ws.Range("E40").Value = "'4.38%"
$ws.Range("D41").Value = "'0.006651"
$ws.Range("E41").Value = "'58.60%"
$ws.Range("D42").Value = "'0.1179"
$ws.Range("E42").Value = "'5.68%"
$ws.Range("D43").Value = "'0.002358"
$ws.Range("E43").Value = "'11.96%"
$ws.Range("D44").Value = "'0.01307"
$ws.Range("E44").Value = "'13.41%"
$ws.Range("D45").Value = "'0.00005248"
$ws.Range("E45").Value = "'2.64%"
$ws.Range("E46").Value = "'-0.08%"
$ws.Range("E47").Value = "'812.21%"
$ws.Range("D48").Value = "'0.02249"
$ws.Range("E48").Value = "'-8.12%"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("E50").Value = "'-0.15%"
